$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Normalize the half-width "例:" separator to full-width "例：" (and tidy
# the ":"-after-letter spacing in the PieceCode legend) in the 備註說明
# (remarks) column for the PieceCode-related field rows.
$ws.Range("G10").Value = "~~~(綁約2年以下或未綁約)~~~`nA:新貸件`nB:新貸件(同押品,數額度之額度一以外)`nC:原額度內—動支件`nD:新增額度—新貸件(指有增加設定抵押權者)`nE:展期件`n~~~(有綁約2年(含)以上)~~~`n1:新貸件`n2:新貸件(同押品,數額度之額度一以外)`n3:原額度內—動支件`n4:新增額度—新貸件(指有增加設定抵押權者)`n5:展期件`n~~~(無關綁約)~~~`n6:原額度內—6個月內動支件(還款後6個月內再動支者)`n7:服務件`n8:特殊件`n9:固特利契轉"
$ws.Range("G11").Value = "0=不計件數`n例： 1=1件 , 2=2件"
$ws.Range("G12").Value = "介紹單位_件數>0時有值`n例： 600000=1件(60萬以上)"
$ws.Range("G13").Value = "輸入比例`n例： 0=無 , 1=全部計算"
$ws.Range("G14").Value = "輸入比例`n例： 0=無 , 0.001=1/1000*撥款金額"
$ws.Range("G15").Value = "例： 0=無限制 , 500000=撥款金額達50萬以上者"
$ws.Range("G16").Value = "例： 0=無 , 10000=以每一萬元計算獎金金額"
$ws.Range("G17").Value = "介紹人_換算業績金額基準=0時必須為0`n例： 0=無 , 35=以每一萬元計算35元業績獎金"
$ws.Range("G18").Value = "例： 0=無 , 10000=以每一萬元計算獎金金額"
$ws.Range("G19").Value = "介紹人_二階(或三階承攬)業務報酬_金額基準=0時必須為0`n例： 0=無 , 12.5=以每一萬元計算12.5元業績獎金"
$ws.Range("G20").Value = "0=不計件數`n例： 1=1件 , 0.1=0.1件"
$ws.Range("G21").Value = "0=無上限`n例： 1=最高1件"
$ws.Range("G22").Value = "例： 600000=1件(60萬以上)`n例： 100000=10萬計0.1件,最高1件"
$ws.Range("G23").Value = "輸入比例`n例： 0=無 , 1=全部計算 , 0.5=撥款金額*1/2"

# Leave the cursor where the author left it when saving.
$ws.Range("H45").Select()
